$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("E2").Value = 7
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 13

# Move the active selection from A2 to B2
$ws.Range("B2").Select()
